$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes current rows 4-15 down to 5-16)
$ws.Range("A4").EntireRow.Insert()

# Copy the date number format from the row below (now row 5, formerly row 4)
$ws.Range("D4").NumberFormat = $ws.Range("D5").NumberFormat

# Fill in the new row 4 data
$ws.Range("A4").Value = 12
$ws.Range("B4").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 45251
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 100112032
$ws.Range("G4").Value = "Zapallo italiano"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("N4").Value = "$/caja 50 unidades"
$ws.Range("O4").Value = "Región de O'Higgins"
$ws.Range("P4").Value = 280
$ws.Range("Q4").Value = 50
$ws.Range("R4").Value = "Hortaliza"
